# Minor edits to Clinical Characteristics table.
$d = $word.ActiveDocument

# wdReplaceAll = 2, wdFindContinue not needed since MatchCase/whole word off;
# these strings are each unique across the document so a plain Find/Replace
# is safe and unambiguous.

# 1) Row header label: "Radical nephrectomy" -> "Type of radical nephrectomy"
$d.Content.Find.Execute("Radical nephrectomy", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Type of radical nephrectomy", 2)

# 2) p-value for the (now renamed) radical nephrectomy row: 0.086 -> 0.14
$d.Content.Find.Execute("0.086", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.14", 2)

# 3) Laparoscopic/Robotic row (under radical nephrectomy) counts
$d.Content.Find.Execute("56 (75%)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "61 (81%)", 2)
$d.Content.Find.Execute("32 (86%)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "33 (89%)", 2)
$d.Content.Find.Execute("24 (63%)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "28 (74%)", 2)

# 4) Remove the "Laparoscopic/Robotic with Adrenalectomy" row entirely -
#    it is table 2 (Clinical Characteristics), row 10.
$table = $d.Tables.Item(2)
$table.Rows.Item(10).Delete()

# 5) Row header label: "Partial nephrectomy" -> "Type of partial nephrectomy"
$d.Content.Find.Execute("Partial nephrectomy", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Type of partial nephrectomy", 2)
